# Generate Report for Handoff
# Replaces the old GUID-based file references (1bfb144a-ca6b-4aa3-affd-e5abd2347b39)
# with the new GUID (ba056d89-61db-4787-85ec-ff51c52bd823) across the Overview,
# zh-cn and de-de sheets, and bumps the associated handoff/handback timestamps.

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

$oldGuid = "1bfb144a-ca6b-4aa3-affd-e5abd2347b39"
$newGuid = "ba056d89-61db-4787-85ec-ff51c52bd823"

# The external hyperlink target (relationship) is left untouched by the change -
# only the displayed text changes - so keep using the original URL (still
# referencing the old GUID) when re-creating the hyperlinks below.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7fc1b104baee9f839dd4cbfc4a2d31f557dac787/e2e/$oldGuid.md"

### Overview sheet ###
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-20 04:59:58"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    $hyperlinkUrl,
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "e2e\$newGuid.md"
) | Out-Null

### zh-cn sheet ###
$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.346578315c49711ff87c52feab484250854e17e6.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-20 04:59:54"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    $hyperlinkUrl,
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid.md"
) | Out-Null

### de-de sheet ###
$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.346578315c49711ff87c52feab484250854e17e6.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-20 04:59:58"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    $hyperlinkUrl,
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid.md"
) | Out-Null
